$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.187.13"
$ws.Range("E2").Value = "  +5.05%  "
$ws.Range("D3").Value = "3.543.39"
$ws.Range("E3").Value = "  +5.93%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'189.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.58%  "
$ws.Range("D6").Value = "'559.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.39%  "
$ws.Range("D7").Value = "3.535.44"
$ws.Range("E7").Value = "  +5.75%  "
$ws.Range("D8").Value = "'0.616"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.02%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "'0.631"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.99%  "
$ws.Range("D11").Value = "'0.153"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +13.04%  "
$ws.Range("D12").Value = "'54.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("D13").Value = "'0.0000272"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.99%  "
$ws.Range("D14").Value = "'9.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("D15").Value = "4.102.95"
$ws.Range("E15").Value = "  +5.63%  "
$ws.Range("D16").Value = "3.539.88"
$ws.Range("E16").Value = "  +5.89%  "
$ws.Range("D17").Value = "'0.121"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.43%  "
$ws.Range("D18").Value = "67.202.27"
$ws.Range("E18").Value = "  +4.97%  "
$ws.Range("D19").Value = "'18.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.95%  "
$ws.Range("D20").Value = "'12.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.04%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.86%  "
$ws.Range("D22").Value = "'429.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +15.14%  "
$ws.Range("D23").Value = "'4.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.90%  "
$ws.Range("D24").Value = "'85.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.18%  "
$ws.Range("E25").Value = "  -1.91%  "
$ws.Range("D26").Value = "'11.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.90%  "
$ws.Range("E27").Value = "  +6.92%  "
$ws.Range("D28").Value = "'12.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.83%  "
$ws.Range("D29").Value = "'6.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("D30").Value = "'9.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.54%  "
$ws.Range("D31").Value = "'30.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.03%  "
$ws.Range("D32").Value = "'643.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("D33").Value = "'6.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D34").Value = "'11.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.84%  "
$ws.Range("E35").Value = "  +3.77%  "
$ws.Range("D36").Value = "'60.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.52%  "
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "'38.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.91%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0813"
$ws.Range("E38").Value = "  +11.22%  "
$ws.Range("D39").Value = "'0.146"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +16.29%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("D41").Value = "'0.389"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("D42").Value = "'3.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +13.41%  "
$ws.Range("D43").Value = "3.116.12"
$ws.Range("E43").Value = "  +6.12%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "'2.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("D46").Value = "'2.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.71%  "
$ws.Range("D47").Value = "'3.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.71%  "
$ws.Range("D48").Value = "'0.0419"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.35%  "
$ws.Range("D49").Value = "'2.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.02%  "
$ws.Range("D50").Value = "'0.131"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.56%  "
$ws.Range("D51").Value = "'141.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.91%  "
